# Apply the "additional scraping" edit:
#  1. Insert a new "Player Info" sheet as the first sheet with basic player
#     details (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  2. On "ODI Batting" and "ODI Bowling", rename the MATCH_CARD_LINK column
#     to MATCH_CODE and replace the full scorecard URL values with just the
#     trailing numeric match code (kept as text, like the rest of the sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored as text
# (mirrors the rest of the workbook, where even numeric-looking values
# like "1", "11" etc. are stored as text instead of numbers).
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# -----------------------------------------------------------------
# 1. Create the new "Player Info" sheet as the first sheet
# -----------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $playerInfo.Cells.Item(1, $i + 1)
    Set-TextValue $cell $headers[$i]
}

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$playerRow = @("4788", "Dylan Evers Budge", "Right Handed", "Right Arm Medium")
for ($i = 0; $i -lt $playerRow.Length; $i++) {
    $cell = $playerInfo.Cells.Item(2, $i + 1)
    Set-TextValue $cell $playerRow[$i]
}

# -----------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE and replace URLs with codes
#    on both the "ODI Batting" and "ODI Bowling" sheets.
# -----------------------------------------------------------------
function Update-MatchCodeColumn($sheet, $col) {
    $headerCell = $sheet.Cells.Item(1, $col)
    Set-TextValue $headerCell "MATCH_CODE"
    # Re-apply the bold / bordered / centered header look that the other
    # header cells on the row use (Set-TextValue resets to "Normal" style).
    $headerCell.Font.Bold = $true
    $headerCell.HorizontalAlignment = -4108
    $headerCell.VerticalAlignment = -4160
    $headerCell.Borders.LineStyle = 1
    $headerCell.Borders.Weight = 2

    $usedRange = $sheet.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $cell = $sheet.Cells.Item($r, $col)
        $val = $cell.Value()
        if ($val -ne $null -and $val -ne "") {
            if ($val -match "MatchCode=(\d+)") {
                Set-TextValue $cell $matches[1]
            }
        }
    }
}

$odiBatting = $wb.Worksheets.Item("ODI Batting")
Update-MatchCodeColumn $odiBatting 4

$odiBowling = $wb.Worksheets.Item("ODI Bowling")
Update-MatchCodeColumn $odiBowling 2
